$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: "Chemistry of Materials" -> "Materials chemistry" (B4/C4) ---
$ws.Range("B4").Value = "Materials chemistry"
$ws.Range("C4").Value = "Materials chemistry"

# --- Row 8 & Row 13: "01/01/2012" -> "01/01/2023" (B8/C8 and B13/C13 share this string) ---
# Use a text formula then collapse it to a static value so Excel does not
# auto-convert the text into a date serial number (which would also change
# the cell's number format / style).
$ws.Range("B8").Formula = "=""01/01/2023"""
$ws.Range("B8").Copy()
$ws.Range("B8").PasteSpecial(-4163)

$ws.Range("C8").Formula = "=""01/01/2023"""
$ws.Range("C8").Copy()
$ws.Range("C8").PasteSpecial(-4163)

$ws.Range("B13").Formula = "=""01/01/2023"""
$ws.Range("B13").Copy()
$ws.Range("B13").PasteSpecial(-4163)

$ws.Range("C13").Formula = "=""01/01/2023"""
$ws.Range("C13").Copy()
$ws.Range("C13").PasteSpecial(-4163)

# --- Row 11: new Objectives text in B11/C11 (cells were previously empty) ---
# Copy the formatting from the row above first so the new cells pick up the
# same column styles (s=2 for B, s=3 for C) instead of inheriting column A's style.
$ws.Range("B10").Copy()
$ws.Range("B11").PasteSpecial(-4122)
$ws.Range("C10").Copy()
$ws.Range("C11").PasteSpecial(-4122)
$ws.Range("B11").Value = "Provide the student with the main types of organic and inorganic synthesis of materials as well as presenting the main analytical techniques for material characterization."
$ws.Range("C11").Value = "Provide the student with the main types of organic and inorganic synthesis of materials as well as presenting the main analytical techniques for material characterization."

# --- Row 14: new Short syllabus text in B14/C14 (cells were previously empty) ---
$ws.Range("B15").Copy()
$ws.Range("B14").PasteSpecial(-4122)
$ws.Range("C15").Copy()
$ws.Range("C14").PasteSpecial(-4122)
$ws.Range("B14").Value = "Introduction to the chemistry of materials and its association with the synthesis of new materials. The modern view of the atom and chemical bonds. Crystal structure and crystallographic characterization techniques. Epitaxial thin films and films in general and their impact on modern technology. Amorphous materials, synthesis and applications. Synthesis of materials and chemical transformations. Processes and Techniques of crystal growth in general. Conducting polymers and their applications in modern technology."
$ws.Range("C14").Value = "Introduction to the chemistry of materials and its association with the synthesis of new materials. The modern view of the atom and chemical bonds. Crystal structure and crystallographic characterization techniques. Epitaxial thin films and films in general and their impact on modern technology. Amorphous materials, synthesis and applications. Synthesis of materials and chemical transformations. Processes and Techniques of crystal growth in general. Conducting polymers and their applications in modern technology."

# --- Row 16: new Syllabus text in B16/C16 (cells were previously empty) ---
$ws.Range("B15").Copy()
$ws.Range("B16").PasteSpecial(-4122)
$ws.Range("C15").Copy()
$ws.Range("C16").PasteSpecial(-4122)
$ws.Range("B16").Value = "Materials chemistry: definition; role of chemistry in materials science; fundamentals.Atomistics and the modern view of the atom with quantum foundations.Types of chemical bonds: van der Waals forces, Lennard-Jones potential, covalent bonding, coordination bonds, ionic bonds and metallic bonds.Polycrystalline and monocrystalline materials. The crystallographic order and crystallographic and microscopic characterization techniques. The importance of single crystals in electronic applications. High quality crystal growth techniques such as: flow method, Czochralski method, Brigdmann method, vapor transport method and modified isothermal vapor transport growth method. Amorphous materials and their importance for modern technology. Concepts and techniques for growing amorphous materials. Epitaxial thin films, growth techniques such as: chemical vapor, sputtering, laser ablation and MBE. Thin films grown by electrolysis for protective coating, concepts and applications. Synthesis of conductive polymers, concepts and applications as electronic devices."
$ws.Range("C16").Value = "Materials chemistry: definition; role of chemistry in materials science; fundamentals.Atomistics and the modern view of the atom with quantum foundations.Types of chemical bonds: van der Waals forces, Lennard-Jones potential, covalent bonding, coordination bonds, ionic bonds and metallic bonds.Polycrystalline and monocrystalline materials. The crystallographic order and crystallographic and microscopic characterization techniques. The importance of single crystals in electronic applications. High quality crystal growth techniques such as: flow method, Czochralski method, Brigdmann method, vapor transport method and modified isothermal vapor transport growth method. Amorphous materials and their importance for modern technology. Concepts and techniques for growing amorphous materials. Epitaxial thin films, growth techniques such as: chemical vapor, sputtering, laser ablation and MBE. Thin films grown by electrolysis for protective coating, concepts and applications. Synthesis of conductive polymers, concepts and applications as electronic devices."

# --- Row 20: Norma de recuperação value update (B20/C20) ---
$ws.Range("B20").Value = "Média simples de duas provas escritas,  Conceito Final = (P1 + P2)/2"
$ws.Range("C20").Value = "Média simples de duas provas escritas,  Conceito Final = (P1 + P2)/2"

# --- Row 21: Bibliografia value update (B21/C21) ---
$ws.Range("B21").Value = "Aplicação de duas provas escritas dentro do prazo regimental antes do início do próximo semestre letivo."
$ws.Range("C21").Value = "Aplicação de duas provas escritas dentro do prazo regimental antes do início do próximo semestre letivo."
